$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (header "Förändrad") rows 2 through 97 all hold the same
# serial date value (45188) which needs to be bumped to 45189 (one day later).
$ws.Range("C2:C97").Value = 45189
